$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values with refined/rounded data (custom accuracy)
$ws.Cells.Item(5,1).Value = 40751.4027662037
$ws.Cells.Item(5,2).Value = 5.82
$ws.Cells.Item(5,3).Value = 4.07
$ws.Cells.Item(5,4).Value = 0.06
$ws.Cells.Item(5,5).Value = 10.83
$ws.Cells.Item(5,6).Value = 9.32
$ws.Cells.Item(5,7).Value = 4.44
$ws.Cells.Item(5,8).Value = 21.32
$ws.Cells.Item(5,9).Value = 6.02
$ws.Cells.Item(5,10).Value = 2.91
$ws.Cells.Item(5,11).Value = 4.69
$ws.Cells.Item(5,12).Value = 4.47
$ws.Cells.Item(5,13).Value = 4.49
$ws.Cells.Item(5,14).Value = 1.34
$ws.Cells.Item(5,15).Value = 3.9
$ws.Cells.Item(5,16).Value = 6.25
$ws.Cells.Item(5,17).Value = 3.28
$ws.Cells.Item(5,18).Value = 0.28
$ws.Cells.Item(5,19).Value = 0.07000000000000001
$ws.Cells.Item(5,20).Value = 55.48
$ws.Cells.Item(5,21).Value = 11.82
$ws.Cells.Item(5,22).Value = 3.98
$ws.Cells.Item(5,23).Value = 8.1
$ws.Cells.Item(5,24).Value = 3.98
$ws.Cells.Item(5,25).Value = 0.54
$ws.Cells.Item(5,26).Value = 9.949999999999999
$ws.Cells.Item(5,27).Value = 3.35
$ws.Cells.Item(5,28).Value = 2.79
$ws.Cells.Item(5,29).Value = 3.33
$ws.Cells.Item(5,30).Value = 5.17
$ws.Cells.Item(5,31).Value = 0.52
$ws.Cells.Item(5,32).Value = 19.49
$ws.Cells.Item(5,33).Value = 2.02
$ws.Cells.Item(5,34).Value = 4.6

# Remove row 6 (dataset trimmed)
$ws.Rows.Item(6).Delete()
